# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets
# F3: 1432 -> 1435
# F5: 13   -> 14
# F8: 209  -> 214

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 1435
    $ws.Range("F5").Value = 14
    $ws.Range("F8").Value = 214
}
